$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "A 13467-2023"
$ws.Range("B2").Value = 45005
$ws.Range("C2").Value = 46079
$ws.Range("G2").Value = 2.3
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 2
$ws.Range("J2").Value = 3
$ws.Range("O2").Value = 3
$ws.Range("R2").Value = "Gulsparv`r`nHypoxylon petriniae`r`nKråka`r`nGrå skärelav`r`nGulnål"
$ws.Range("S2").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1283/artfynd/A 13467-2023 artfynd.xlsx`", `"A 13467-2023`")"
$ws.Range("T2").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1283/kartor/A 13467-2023 karta.png`", `"A 13467-2023`")"
$ws.Range("V2").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1283/klagomål/A 13467-2023 FSC-klagomål.docx`", `"A 13467-2023`")"
$ws.Range("W2").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1283/klagomålsmail/A 13467-2023 FSC-klagomål mail.docx`", `"A 13467-2023`")"
$ws.Range("X2").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1283/tillsyn/A 13467-2023 tillsynsbegäran.docx`", `"A 13467-2023`")"
$ws.Range("Y2").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1283/tillsynsmail/A 13467-2023 tillsynsbegäran mail.docx`", `"A 13467-2023`")"
$ws.Range("Z2").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1283/fåglar/A 13467-2023 prioriterade fågelarter.docx`", `"A 13467-2023`")"

# Row 3
$ws.Range("A3").Value = "A 45325-2025"
$ws.Range("B3").Value = 45922
$ws.Range("C3").Value = 46079
$ws.Range("G3").Value = 1.6
$ws.Range("H3").Value = 4
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 1
$ws.Range("O3").Value = 1
$ws.Range("R3").Value = "Nordlig buksimmare`r`nStörre vattensalamander`r`nÅkergroda`r`nMindre vattensalamander`r`nVanlig groda"
$ws.Range("S3").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1283/artfynd/A 45325-2025 artfynd.xlsx`", `"A 45325-2025`")"
$ws.Range("T3").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1283/kartor/A 45325-2025 karta.png`", `"A 45325-2025`")"
$ws.Range("V3").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1283/klagomål/A 45325-2025 FSC-klagomål.docx`", `"A 45325-2025`")"
$ws.Range("W3").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1283/klagomålsmail/A 45325-2025 FSC-klagomål mail.docx`", `"A 45325-2025`")"
$ws.Range("X3").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1283/tillsyn/A 45325-2025 tillsynsbegäran.docx`", `"A 45325-2025`")"
$ws.Range("Y3").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1283/tillsynsmail/A 45325-2025 tillsynsbegäran mail.docx`", `"A 45325-2025`")"
$ws.Range("Z3").ClearContents()

# Row 4
$ws.Range("C4").Value = 46079

# Row 5
$ws.Range("C5").Value = 46079

# Row 6
$ws.Range("A6").Value = "A 12651-2022"
$ws.Range("B6").Value = 44641
$ws.Range("C6").Value = 46079
$ws.Range("G6").Value = 3.2

# Row 7
$ws.Range("A7").Value = "A 5792-2024"
$ws.Range("B7").Value = 45335
$ws.Range("C7").Value = 46079
$ws.Range("G7").Value = 5.6

# Row 8
$ws.Range("A8").Value = "A 2593-2024"
$ws.Range("B8").Value = 45313.69204861111
$ws.Range("C8").Value = 46079
$ws.Range("G8").Value = 2.3

# Row 9
$ws.Range("A9").Value = "A 7333-2025"
$ws.Range("B9").Value = 45703.35899305555
$ws.Range("C9").Value = 46079
$ws.Range("G9").Value = 0.9

# Row 10
$ws.Range("A10").Value = "A 35642-2023"
$ws.Range("B10").Value = 45147
$ws.Range("C10").Value = 46079
$ws.Range("G10").Value = 1.2

# Row 11
$ws.Range("A11").Value = "A 28288-2023"
$ws.Range("B11").Value = 45099.6349537037
$ws.Range("C11").Value = 46079
$ws.Range("G11").Value = 0.5

# Row 12
$ws.Range("A12").Value = "A 13651-2023"
$ws.Range("B12").Value = 45006
$ws.Range("C12").Value = 46079
$ws.Range("G12").Value = 2.2

# Row 13
$ws.Range("A13").Value = "A 8194-2025"
$ws.Range("B13").Value = 45708
$ws.Range("C13").Value = 46079
$ws.Range("G13").Value = 1.9

# Row 14
$ws.Range("A14").Value = "A 50997-2025"
$ws.Range("B14").Value = 45946
$ws.Range("C14").Value = 46079
$ws.Range("G14").Value = 1.5

# Row 15
$ws.Range("A15").Value = "A 7814-2026"
$ws.Range("B15").Value = 46062.61388888889
$ws.Range("C15").Value = 46079
$ws.Range("G15").Value = 1.1

# Row 16
$ws.Range("A16").Value = "A 7827-2026"
$ws.Range("B16").Value = 46062.63958333333
$ws.Range("C16").Value = 46079
$ws.Range("G16").Value = 2.1
